# Actualización automática 2025-08-18 16:20:08
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("P13").Value = 1628.61
$ws1.Range("M14").Value = 782.58
$ws1.Range("D25").Value = 457.92
$ws1.Range("M25").Value = 489.11
$ws1.Range("D38").Value = 1408.32
$ws1.Range("D57").Value = "6 de 55"
$ws1.Range("M57").Value = "13 de 55"
$ws1.Range("P57").Value = "3 de 55"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F13").Value = 1628.61
$ws2.Range("F14").Value = 1840.95
$ws2.Range("F25").Value = 3996.81
$ws2.Range("F38").Value = 1408.32
$ws2.Range("F57").Value = 38762.57

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 9279.35
$ws3.Range("E3").Value = 18177.6576
$ws3.Range("F3").Value = 0.3379592610813132

$ws3.Range("D10").Value = 1668.94
$ws3.Range("E10").Value = -368.4400000000001
$ws3.Range("F10").Value = 1.283306420607459

$ws3.Range("D16").Value = 19100.62
$ws3.Range("E16").Value = 36959.08
$ws3.Range("F16").Value = 0.3407192689222382

$ws3.Range("D19").Value = 38762.57000000001
$ws3.Range("E19").Value = 78677.12064517915
$ws3.Range("F19").Value = 0.330063624887377
